# Update the HIGH/LOW/CLOSE/LTP/VOL/9:25 CLOSE figures on the active
# sheet (Sheet1) for rows 2-17 to reflect the latest market data.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$data = @(
    @{ Row = 2; B = 1445.5; C = 1397.9; D = 1430.35; E = 1433.75; F = 101; G = 1398.75 },
    @{ Row = 3; B = 1267; C = 1242; D = 1259; E = 1260.2; F = 20; G = 1245.6 },
    @{ Row = 4; B = 50175; C = 49750; D = 50035; E = 50042.4; F = 17; G = 49811.85 },
    @{ Row = 5; B = 121.4; C = 119.5; D = 121; E = 121.1; F = 463; G = 119.65 },
    @{ Row = 6; B = 883.9; C = 872.45; D = 877.8; E = 878.65; F = 54; G = 873.6 },
    @{ Row = 7; B = 687.5; C = 679.25; D = 683.15; E = 683.5; F = 100; G = 680.2 },
    @{ Row = 8; B = 1111.6; C = 1099.1; D = 1102.95; E = 1103.35; F = 232; G = 1100.6 },
    @{ Row = 9; B = 1061.9; C = 1036.7; D = 1054.25; E = 1055; F = 29; G = 1039 },
    @{ Row = 10; B = 23489.85; C = 23350.8; D = 23457.1; E = 23466.65; F = 45; G = 23376.3 },
    @{ Row = 11; B = 2959.7; C = 2919.1; D = 2953; E = 2955.4; F = 71; G = 2922.25 },
    @{ Row = 12; B = 848.4; C = 838.95; D = 840.45; E = 840.5; F = 211; G = 843.75 },
    @{ Row = 13; B = 1124.8; C = 1113.35; D = 1115.95; E = 1114.9; F = 16; G = 1120.6 },
    @{ Row = 14; B = 998.3; C = 984.4; D = 992.85; E = 994.2; F = 183; G = 985.4 },
    @{ Row = 15; B = 183.7; C = 181.95; D = 183.25; E = 183.3; F = 360; G = 182.2 },
    @{ Row = 16; B = 3870; C = 3828.95; D = 3834; E = 3832.55; F = 25; G = 3866.8 },
    @{ Row = 17; B = 3525; C = 3469.3; D = 3521; E = 3519.05; F = 26; G = 3482.3 }
)

foreach ($item in $data) {
    $r = $item.Row
    $ws.Cells.Item($r, 2).Value = $item.B
    $ws.Cells.Item($r, 3).Value = $item.C
    $ws.Cells.Item($r, 4).Value = $item.D
    $ws.Cells.Item($r, 5).Value = $item.E
    $ws.Cells.Item($r, 6).Value = $item.F
    $ws.Cells.Item($r, 7).Value = $item.G
}
